$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Header text updates (shared-string rich-text cells)
#    "Volume 30   Number  35" -> "...Number  36"
#    "Report Covering the Week  8/28/2023  Through  9/3/2023"
#      -> "...9/4/2023  Through  9/10/2023"
# ---------------------------------------------------------------
$volChars = $ws.Range("A8").Characters(21, 2)
$volChars.Text = "36"

$date1 = $ws.Range("C9").Characters(27, 9)
$date1.Text = "9/4/2023"

$date2 = $ws.Range("C9").Characters(46, 8)
$date2.Text = "9/10/2023"

# ---------------------------------------------------------------
# 2. Cells that flip from numeric -> text ("0" / "***.*") or
#    text -> numeric. Copy formatting + value from a stable
#    same-style source cell so the underlying style index and
#    storage type line up exactly, then overwrite where needed.
# ---------------------------------------------------------------

# C14: 1 -> "0" (shared string, style 14)
$ws.Range("D14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("C14").PasteSpecial(-4163)

# G14: 1 -> "0" (shared string, style 14)
$ws.Range("D14").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("G14").PasteSpecial(-4163)

# H14: 200 -> "***.*" (shared string, style 14)
$ws.Range("E14").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("H14").PasteSpecial(-4163)

# C15: 2 -> "0" (shared string, style 14)
$ws.Range("D14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("C15").PasteSpecial(-4163)

# C26: 2 -> "0" (shared string, style 14)
$ws.Range("D14").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("C26").PasteSpecial(-4163)

# D22: "0" (text) -> 4 (numeric, style 15)
$ws.Range("C16").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 4

# E22: "***.*" (text) -> -75 (numeric, style 16)
$ws.Range("K14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = -75

$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# 3. Plain numeric value updates (same type / style throughout)
# ---------------------------------------------------------------

# Row 14 (Murder)
$ws.Range("F14").Value = 4
$ws.Range("I14").Value = 7
$ws.Range("K14").Value = -41.666666666666
$ws.Range("L14").Value = 16.666666666666
$ws.Range("M14").Value = -12.5
$ws.Range("N14").Value = -73.076923076923

# Row 15 (Rape)
$ws.Range("E15").Value = -100
$ws.Range("J15").Value = 30
$ws.Range("K15").Value = -16.666666666666
$ws.Range("M15").Value = 8.695652173913
$ws.Range("N15").Value = -43.181818181818

# Row 16 (Robbery)
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = 10
$ws.Range("F16").Value = 39
$ws.Range("G16").Value = 29
$ws.Range("H16").Value = 34.482758620689
$ws.Range("I16").Value = 321
$ws.Range("J16").Value = 294
$ws.Range("K16").Value = 9.183673469387
$ws.Range("L16").Value = 27.888446215139
$ws.Range("M16").Value = -5.865102639296
$ws.Range("N16").Value = -71.466666666666

# Row 17 (Fel. Assault)
$ws.Range("C17").Value = 12
$ws.Range("E17").Value = -14.285714285714
$ws.Range("F17").Value = 52
$ws.Range("G17").Value = 43
$ws.Range("H17").Value = 20.930232558139
$ws.Range("I17").Value = 470
$ws.Range("J17").Value = 433
$ws.Range("K17").Value = 8.545034642032
$ws.Range("L17").Value = 28.767123287671
$ws.Range("M17").Value = 45.962732919254
$ws.Range("N17").Value = 3.070175438596

# Row 18 (Burglary)
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 160
$ws.Range("J18").Value = 144
$ws.Range("K18").Value = 11.111111111111
$ws.Range("L18").Value = 31.147540983606
$ws.Range("M18").Value = -40.074906367041
$ws.Range("N18").Value = -89.974937343358

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value = 19
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = 35.714285714285
$ws.Range("F19").Value = 63
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = 36.956521739130
$ws.Range("I19").Value = 507
$ws.Range("J19").Value = 460
$ws.Range("K19").Value = 10.217391304347
$ws.Range("L19").Value = 10.217391304347
$ws.Range("M19").Value = 32.722513089005
$ws.Range("N19").Value = -19.138755980861

# Row 20 (G.L.A.)
$ws.Range("C20").Value = 11
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 120
$ws.Range("F20").Value = 38
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = 111.111111111111
$ws.Range("I20").Value = 245
$ws.Range("J20").Value = 182
$ws.Range("K20").Value = 34.615384615384
$ws.Range("L20").Value = 62.251655629139
$ws.Range("M20").Value = 100.819672131148
$ws.Range("N20").Value = -76.532567049808

# Row 21 (TOTAL)
$ws.Range("C21").Value = 57
$ws.Range("D21").Value = 48
$ws.Range("E21").Value = 18.75
$ws.Range("F21").Value = 216
$ws.Range("G21").Value = 155
$ws.Range("H21").Value = 39.354838709677
$ws.Range("I21").Value = 1735
$ws.Range("J21").Value = 1555
$ws.Range("K21").Value = 11.575562700964
$ws.Range("L21").Value = 26.181818181818
$ws.Range("M21").Value = 18.430034129692
$ws.Range("N21").Value = -64.721431476209

# Row 22 (Transit)
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = -57.142857142857
$ws.Range("I22").Value = 31
$ws.Range("J22").Value = 34
$ws.Range("K22").Value = -8.823529411764
$ws.Range("L22").Value = 158.333333333333
$ws.Range("M22").Value = 72.222222222222

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 53
$ws.Range("D24").Value = 44
$ws.Range("E24").Value = 20.454545454545
$ws.Range("G24").Value = 180
$ws.Range("H24").Value = 26.111111111111
$ws.Range("I24").Value = 1790
$ws.Range("J24").Value = 2021
$ws.Range("K24").Value = -11.429985155863
$ws.Range("L24").Value = 80.990899898887
$ws.Range("M24").Value = 101.349831271091

# Row 25 (Misd. Assault)
$ws.Range("C25").Value = 25
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = 47.058823529411
$ws.Range("F25").Value = 66
$ws.Range("G25").Value = 61
$ws.Range("H25").Value = 8.196721311475
$ws.Range("I25").Value = 590
$ws.Range("J25").Value = 551
$ws.Range("K25").Value = 7.078039927404
$ws.Range("L25").Value = 10.486891385767
$ws.Range("M25").Value = -6.793048973143

# Row 26 (UCR Rape*)
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = 175
$ws.Range("I26").Value = 43
$ws.Range("J26").Value = 60
$ws.Range("K26").Value = -28.333333333333
$ws.Range("L26").Value = 53.571428571428

# Row 27 (Other Sex Crimes)
$ws.Range("C27").Value = 4
$ws.Range("E27").Value = 33.333333333333
$ws.Range("F27").Value = 11
$ws.Range("G27").Value = 12
$ws.Range("H27").Value = -8.333333333333
$ws.Range("I27").Value = 71
$ws.Range("J27").Value = 56
$ws.Range("K27").Value = 26.785714285714
$ws.Range("L27").Value = -14.457831325301

# Row 28 (Shooting Vic.)
$ws.Range("D28").Value = 2
$ws.Range("J28").Value = 28
$ws.Range("K28").Value = -50
$ws.Range("M28").Value = -51.724137931034
$ws.Range("N28").Value = -75

# Row 29 (Shooting Inc.)
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("J29").Value = 24
$ws.Range("K29").Value = -54.166666666666
$ws.Range("M29").Value = -50
$ws.Range("N29").Value = -78
